$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(9, 8).Value = 229.6  # H9: 242 -> 229.6
$ws.Cells.Item(9, 10).Value = 239.5  # J9: 259.33334 -> 239.5
$ws.Cells.Item(9, 12).Value = 239.5  # L9: 259.33334 -> 239.5
$ws.Cells.Item(9, 14).Value = -577.5  # N9: -597.33334 -> -577.5
$ws.Cells.Item(59, 8).Value = 1666.6666  # H59: 1800 -> 1666.6666
$ws.Cells.Item(59, 9).Value = 1250  # I59: 1500 -> 1250
$ws.Cells.Item(59, 10).Value = 2500  # J59: 1950 -> 2500
$ws.Cells.Item(59, 11).Value = 3750  # K59: 4500 -> 3750
$ws.Cells.Item(59, 12).Value = 7500  # L59: 5850 -> 7500
$ws.Cells.Item(59, 13).Value = -3193  # M59: -3943 -> -3193
$ws.Cells.Item(59, 14).Value = -8614  # N59: -6964 -> -8614
$ws.Cells.Item(62, 8).Value = 4989.4287  # H62: 4378.4546 -> 4989.4287
$ws.Cells.Item(62, 9).Value = 4672.7896  # I62: 4805.1665 -> 4672.7896
$ws.Cells.Item(62, 10).Value = 7997.5  # J62: 3866.4 -> 7997.5
$ws.Cells.Item(62, 11).Value = 4672.7896  # K62: 4805.1665 -> 4672.7896
$ws.Cells.Item(62, 12).Value = 7997.5  # L62: 3866.4 -> 7997.5
$ws.Cells.Item(62, 13).Value = -4048.7896  # M62: -4181.1665 -> -4048.7896
$ws.Cells.Item(62, 14).Value = -9245.5  # N62: -5114.4 -> -9245.5
$ws.Cells.Item(65, 8).Value = 4989.4287  # H65: 4378.4546 -> 4989.4287
$ws.Cells.Item(65, 9).Value = 4672.7896  # I65: 4805.1665 -> 4672.7896
$ws.Cells.Item(65, 10).Value = 7997.5  # J65: 3866.4 -> 7997.5
$ws.Cells.Item(65, 11).Value = 23363.948  # K65: 24025.8325 -> 23363.948
$ws.Cells.Item(65, 12).Value = 39987.5  # L65: 19332 -> 39987.5
$ws.Cells.Item(65, 13).Value = -20243.948  # M65: -20905.8325 -> -20243.948
$ws.Cells.Item(65, 14).Value = -46227.5  # N65: -25572 -> -46227.5
$ws.Cells.Item(86, 8).Value = 18511.555  # H86: 17545.5 -> 18511.555
$ws.Cells.Item(86, 9).Value = 22000  # I86: 14400 -> 22000
$ws.Cells.Item(86, 10).Value = 16767.334  # J86: 20691 -> 16767.334
$ws.Cells.Item(86, 11).Value = 22000  # K86: 14400 -> 22000
$ws.Cells.Item(86, 12).Value = 16767.334  # L86: 20691 -> 16767.334
$ws.Cells.Item(86, 13).Value = -20877  # M86: -13277 -> -20877
$ws.Cells.Item(86, 14).Value = -19013.334  # N86: -22937 -> -19013.334
$ws.Cells.Item(89, 8).Value = 18511.555  # H89: 17545.5 -> 18511.555
$ws.Cells.Item(89, 9).Value = 22000  # I89: 14400 -> 22000
$ws.Cells.Item(89, 10).Value = 16767.334  # J89: 20691 -> 16767.334
$ws.Cells.Item(89, 11).Value = 110000  # K89: 72000 -> 110000
$ws.Cells.Item(89, 12).Value = 83836.67  # L89: 103455 -> 83836.67
$ws.Cells.Item(89, 13).Value = -104384  # M89: -66384 -> -104384
$ws.Cells.Item(89, 14).Value = -95068.67  # N89: -114687 -> -95068.67
$ws.Cells.Item(97, 8).Value = 1099  # H97: 1149.5 -> 1099
$ws.Cells.Item(97, 10).Value = 1099  # J97: 1149.5 -> 1099
$ws.Cells.Item(97, 12).Value = 3297  # L97: 3448.5 -> 3297
$ws.Cells.Item(97, 14).Value = -4289  # N97: -4440.5 -> -4289
$ws.Cells.Item(100, 8).Value = 2471.1428  # H100: 3699.75 -> 2471.1428
$ws.Cells.Item(100, 9).Value = 2633  # I100: 4266.3335 -> 2633
$ws.Cells.Item(100, 10).Value = 1500  # J100: 2000 -> 1500
$ws.Cells.Item(100, 11).Value = 2633  # K100: 4266.3335 -> 2633
$ws.Cells.Item(100, 12).Value = 1500  # L100: 2000 -> 1500
$ws.Cells.Item(100, 13).Value = -2092  # M100: -3725.3335 -> -2092
$ws.Cells.Item(100, 14).Value = -2582  # N100: -3082 -> -2582
$ws.Cells.Item(106, 8).Value = 3649.75  # H106: 3319.8 -> 3649.75
$ws.Cells.Item(106, 10).Value = 4000  # J106: 3000 -> 4000
$ws.Cells.Item(106, 12).Value = 4000  # L106: 3000 -> 4000
$ws.Cells.Item(106, 14).Value = -5262  # N106: -4262 -> -5262
$ws.Cells.Item(118, 8).Value = 2600  # H118: 2688.5 -> 2600
$ws.Cells.Item(118, 10).Value = 0  # J118: 2777 -> 0
$ws.Cells.Item(118, 12).Value = 0  # L118: 8331 -> 0
$ws.Cells.Item(118, 14).ClearContents()  # N118: -11645 -> (removed)
$ws.Cells.Item(125, 8).Value = 2650.6667  # H125: 2349.4285 -> 2650.6667
$ws.Cells.Item(125, 9).Value = 780.8  # I125: 741 -> 780.8
$ws.Cells.Item(125, 11).Value = 7027.2  # K125: 6669 -> 7027.2
$ws.Cells.Item(125, 13).Value = -4567.2  # M125: -4209 -> -4567.2
$ws.Cells.Item(132, 8).Value = 5086  # H132: 5132.76 -> 5086
$ws.Cells.Item(132, 9).Value = 4723.048  # I132: 4778.7144 -> 4723.048
$ws.Cells.Item(132, 11).Value = 14169.144  # K132: 14336.1432 -> 14169.144
$ws.Cells.Item(132, 13).Value = -11639.144  # M132: -11806.1432 -> -11639.144

# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = 7516.6523  # H32: 7090.1836 -> 7516.6523
$ws.Cells.Item(32, 9).Value = 7516.6523  # I32: 7090.1836 -> 7516.6523
$ws.Cells.Item(32, 11).Value = 7516.6523  # K32: 7090.1836 -> 7516.6523
$ws.Cells.Item(32, 13).Value = -7229.6523  # M32: -6803.1836 -> -7229.6523
$ws.Cells.Item(36, 8).Value = 3710.4  # H36: 2638 -> 3710.4
$ws.Cells.Item(36, 10).Value = 8000  # J36: 0 -> 8000
$ws.Cells.Item(36, 12).Value = 8000  # L36: 0 -> 8000
$ws.Cells.Item(36, 14).Value = -8692  # N36: None -> -8692
$ws.Cells.Item(44, 8).Value = 34998  # H44: 34997.5 -> 34998
$ws.Cells.Item(44, 10).Value = 34998  # J44: 34997.5 -> 34998
$ws.Cells.Item(44, 12).Value = 34998  # L44: 34997.5 -> 34998
$ws.Cells.Item(44, 14).Value = -35974  # N44: -35973.5 -> -35974
$ws.Cells.Item(61, 8).Value = 2333.3333  # H61: 3000 -> 2333.3333
$ws.Cells.Item(61, 9).Value = 2333.3333  # I61: 3000 -> 2333.3333
$ws.Cells.Item(61, 11).Value = 2333.3333  # K61: 3000 -> 2333.3333
$ws.Cells.Item(61, 13).Value = -2121.3333  # M61: -2788 -> -2121.3333
$ws.Cells.Item(80, 8).Value = 38000  # H80: 38333.332 -> 38000
$ws.Cells.Item(83, 8).Value = 38000  # H83: 38333.332 -> 38000
$ws.Cells.Item(97, 8).Value = 890.5  # H97: 1082.0625 -> 890.5
$ws.Cells.Item(97, 9).Value = 633.6875  # I97: 648 -> 633.6875
$ws.Cells.Item(97, 10).Value = 2945  # J97: 2963 -> 2945
$ws.Cells.Item(97, 11).Value = 633.6875  # K97: 648 -> 633.6875
$ws.Cells.Item(97, 12).Value = 2945  # L97: 2963 -> 2945
$ws.Cells.Item(97, 13).Value = -137.6875  # M97: -152 -> -137.6875
$ws.Cells.Item(97, 14).Value = -3937  # N97: -3955 -> -3937
$ws.Cells.Item(102, 8).Value = 795.6  # H102: 882.25 -> 795.6
$ws.Cells.Item(102, 9).Value = 772.8889  # I102: 865.4286 -> 772.8889
$ws.Cells.Item(102, 11).Value = 772.8889  # K102: 865.4286 -> 772.8889
$ws.Cells.Item(102, 13).Value = 849.1111  # M102: 756.5714 -> 849.1111
$ws.Cells.Item(108, 8).Value = 110000  # H108: 0 -> 110000
$ws.Cells.Item(108, 10).Value = 110000  # J108: 0 -> 110000
$ws.Cells.Item(108, 12).Value = 110000  # L108: 0 -> 110000
$ws.Cells.Item(108, 14).Value = -117680  # N108: None -> -117680
$ws.Cells.Item(122, 8).Value = 2999.5  # H122: 3250 -> 2999.5
$ws.Cells.Item(122, 9).Value = 2999.5  # I122: 3250 -> 2999.5
$ws.Cells.Item(122, 11).Value = 8998.5  # K122: 9750 -> 8998.5
$ws.Cells.Item(122, 13).Value = -6548.5  # M122: -7300 -> -6548.5
$ws.Cells.Item(132, 8).Value = 3352.4546  # H132: 3536.125 -> 3352.4546
$ws.Cells.Item(132, 9).Value = 3553  # I132: 3536.125 -> 3553
$ws.Cells.Item(132, 10).Value = 2450  # J132: 0 -> 2450
$ws.Cells.Item(132, 11).Value = 10659  # K132: 10608.375 -> 10659
$ws.Cells.Item(132, 12).Value = 7350  # L132: 0 -> 7350
$ws.Cells.Item(132, 13).Value = -8129  # M132: -8078.375 -> -8129
$ws.Cells.Item(132, 14).Value = -12410  # N132: None -> -12410
$ws.Cells.Item(136, 8).Value = 2333.3333  # H136: 3000 -> 2333.3333
$ws.Cells.Item(136, 9).Value = 2333.3333  # I136: 3000 -> 2333.3333
$ws.Cells.Item(136, 11).Value = 6999.999899999999  # K136: 9000 -> 6999.999899999999
$ws.Cells.Item(136, 13).Value = -4449.999899999999  # M136: -6450 -> -4449.999899999999
$ws.Cells.Item(141, 8).Value = 0  # H141: 60000 -> 0
$ws.Cells.Item(141, 10).Value = 0  # J141: 60000 -> 0
$ws.Cells.Item(141, 12).Value = 0  # L141: 60000 -> 0
$ws.Cells.Item(141, 14).ClearContents()  # N141: -70360 -> (removed)

# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(64, 8).Value = 1143.7858  # H64: 1217.7858 -> 1143.7858
$ws.Cells.Item(64, 9).Value = 1371.1666  # I64: 1390.4286 -> 1371.1666
$ws.Cells.Item(64, 10).Value = 973.25  # J64: 1045.1428 -> 973.25
$ws.Cells.Item(64, 11).Value = 1371.1666  # K64: 1390.4286 -> 1371.1666
$ws.Cells.Item(64, 12).Value = 973.25  # L64: 1045.1428 -> 973.25
$ws.Cells.Item(64, 13).Value = -1146.1666  # M64: -1165.4286 -> -1146.1666
$ws.Cells.Item(64, 14).Value = -1423.25  # N64: -1495.1428 -> -1423.25
$ws.Cells.Item(67, 8).Value = 1143.7858  # H67: 1217.7858 -> 1143.7858
$ws.Cells.Item(67, 9).Value = 1371.1666  # I67: 1390.4286 -> 1371.1666
$ws.Cells.Item(67, 10).Value = 973.25  # J67: 1045.1428 -> 973.25
$ws.Cells.Item(67, 11).Value = 1371.1666  # K67: 1390.4286 -> 1371.1666
$ws.Cells.Item(67, 12).Value = 973.25  # L67: 1045.1428 -> 973.25
$ws.Cells.Item(67, 13).Value = -591.1666  # M67: -610.4286 -> -591.1666
$ws.Cells.Item(67, 14).Value = -2533.25  # N67: -2605.1428 -> -2533.25
$ws.Cells.Item(94, 8).Value = 1396.1666  # H94: 1241.6428 -> 1396.1666
$ws.Cells.Item(94, 9).Value = 1250.3636  # I94: 1106.3846 -> 1250.3636
$ws.Cells.Item(94, 11).Value = 1250.3636  # K94: 1106.3846 -> 1250.3636
$ws.Cells.Item(94, 13).Value = -799.3635999999999  # M94: -655.3846000000001 -> -799.3635999999999
$ws.Cells.Item(134, 8).Value = 6347.567  # H134: 6949.074 -> 6347.567
$ws.Cells.Item(134, 9).Value = 6623.3335  # I134: 7133.12 -> 6623.3335
$ws.Cells.Item(134, 10).Value = 3865.6667  # J134: 4648.5 -> 3865.6667
$ws.Cells.Item(134, 11).Value = 19870.0005  # K134: 21399.36 -> 19870.0005
$ws.Cells.Item(134, 12).Value = 11597.0001  # L134: 13945.5 -> 11597.0001
$ws.Cells.Item(134, 13).Value = -17335.0005  # M134: -18864.36 -> -17335.0005
$ws.Cells.Item(134, 14).Value = -16667.0001  # N134: -19015.5 -> -16667.0001

# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(14, 8).Value = 0  # H14: 10000 -> 0
$ws.Cells.Item(14, 10).Value = 0  # J14: 10000 -> 0
$ws.Cells.Item(14, 12).Value = 0  # L14: 10000 -> 0
$ws.Cells.Item(14, 14).ClearContents()  # N14: -10340 -> (removed)
$ws.Cells.Item(22, 8).Value = 1003  # H22: 1075.2858 -> 1003
$ws.Cells.Item(22, 9).Value = 825  # I22: 887.5 -> 825
$ws.Cells.Item(22, 10).Value = 1092  # J22: 1150.4 -> 1092
$ws.Cells.Item(22, 11).Value = 825  # K22: 887.5 -> 825
$ws.Cells.Item(22, 12).Value = 1092  # L22: 1150.4 -> 1092
$ws.Cells.Item(22, 13).Value = -475  # M22: -537.5 -> -475
$ws.Cells.Item(22, 14).Value = -1792  # N22: -1850.4 -> -1792
$ws.Cells.Item(25, 8).Value = 0  # H25: 1000 -> 0
$ws.Cells.Item(25, 9).Value = 0  # I25: 1000 -> 0
$ws.Cells.Item(25, 11).Value = 0  # K25: 1000 -> 0
$ws.Cells.Item(25, 13).ClearContents()  # M25: -826 -> (removed)
$ws.Cells.Item(31, 8).Value = 3565.6365  # H31: 3026.7368 -> 3565.6365
$ws.Cells.Item(31, 9).Value = 2501  # I31: 2178.3333 -> 2501
$ws.Cells.Item(31, 10).Value = 5428.75  # J31: 3790.3 -> 5428.75
$ws.Cells.Item(31, 11).Value = 2501  # K31: 2178.3333 -> 2501
$ws.Cells.Item(31, 12).Value = 5428.75  # L31: 3790.3 -> 5428.75
$ws.Cells.Item(31, 13).Value = -2206  # M31: -1883.3333 -> -2206
$ws.Cells.Item(31, 14).Value = -6018.75  # N31: -4380.3 -> -6018.75
$ws.Cells.Item(34, 8).Value = 3565.6365  # H34: 3026.7368 -> 3565.6365
$ws.Cells.Item(34, 9).Value = 2501  # I34: 2178.3333 -> 2501
$ws.Cells.Item(34, 10).Value = 5428.75  # J34: 3790.3 -> 5428.75
$ws.Cells.Item(34, 11).Value = 2501  # K34: 2178.3333 -> 2501
$ws.Cells.Item(34, 12).Value = 5428.75  # L34: 3790.3 -> 5428.75
$ws.Cells.Item(34, 13).Value = -2299  # M34: -1976.3333 -> -2299
$ws.Cells.Item(34, 14).Value = -5832.75  # N34: -4194.3 -> -5832.75
$ws.Cells.Item(41, 8).Value = 15777.777  # H41: 15699.5 -> 15777.777
$ws.Cells.Item(41, 10).Value = 20000  # J41: 19165.834 -> 20000
$ws.Cells.Item(41, 12).Value = 20000  # L41: 19165.834 -> 20000
$ws.Cells.Item(41, 14).Value = -20856  # N41: -20021.834 -> -20856
$ws.Cells.Item(58, 8).Value = 1906.5  # H58: 1483.7 -> 1906.5
$ws.Cells.Item(58, 9).Value = 2092.8  # I58: 1608 -> 2092.8
$ws.Cells.Item(58, 10).Value = 975  # J58: 986.5 -> 975
$ws.Cells.Item(58, 11).Value = 2092.8  # K58: 1608 -> 2092.8
$ws.Cells.Item(58, 12).Value = 975  # L58: 986.5 -> 975
$ws.Cells.Item(58, 13).Value = -1889.8  # M58: -1405 -> -1889.8
$ws.Cells.Item(58, 14).Value = -1381  # N58: -1392.5 -> -1381
$ws.Cells.Item(59, 8).Value = 26569.25  # H59: 28255.4 -> 26569.25
$ws.Cells.Item(93, 8).Value = 17000  # H93: 11400 -> 17000
$ws.Cells.Item(93, 9).Value = 17000  # I93: 8000 -> 17000
$ws.Cells.Item(93, 10).Value = 0  # J93: 25000 -> 0
$ws.Cells.Item(93, 11).Value = 17000  # K93: 8000 -> 17000
$ws.Cells.Item(93, 12).Value = 0  # L93: 25000 -> 0
$ws.Cells.Item(93, 13).Value = -15128  # M93: -6128 -> -15128
$ws.Cells.Item(93, 14).ClearContents()  # N93: -28744 -> (removed)
$ws.Cells.Item(96, 8).Value = 17957  # H96: 27082.666 -> 17957
$ws.Cells.Item(96, 10).Value = 17957  # J96: 27082.666 -> 17957
$ws.Cells.Item(96, 12).Value = 17957  # L96: 27082.666 -> 17957
$ws.Cells.Item(96, 14).Value = -23449  # N96: -32574.666 -> -23449
$ws.Cells.Item(132, 8).Value = 1176.375  # H132: 1333.3334 -> 1176.375
$ws.Cells.Item(132, 9).Value = 918.5  # I132: 1000 -> 918.5
$ws.Cells.Item(132, 10).Value = 1950  # J132: 2000 -> 1950
$ws.Cells.Item(132, 11).Value = 2755.5  # K132: 3000 -> 2755.5
$ws.Cells.Item(132, 12).Value = 5850  # L132: 6000 -> 5850
$ws.Cells.Item(132, 13).Value = -225.5  # M132: -470 -> -225.5
$ws.Cells.Item(132, 14).Value = -10910  # N132: -11060 -> -10910
$ws.Cells.Item(134, 8).Value = 2939.2632  # H134: 3301.2222 -> 2939.2632
$ws.Cells.Item(134, 9).Value = 2952.5  # I134: 3387.7334 -> 2952.5
$ws.Cells.Item(134, 11).Value = 8857.5  # K134: 10163.2002 -> 8857.5
$ws.Cells.Item(134, 13).Value = -6322.5  # M134: -7628.200199999999 -> -6322.5
$ws.Cells.Item(136, 8).Value = 1906.5  # H136: 1483.7 -> 1906.5
$ws.Cells.Item(136, 9).Value = 2092.8  # I136: 1608 -> 2092.8
$ws.Cells.Item(136, 10).Value = 975  # J136: 986.5 -> 975
$ws.Cells.Item(136, 11).Value = 6278.400000000001  # K136: 4824 -> 6278.400000000001
$ws.Cells.Item(136, 12).Value = 2925  # L136: 2959.5 -> 2925
$ws.Cells.Item(136, 13).Value = -3728.400000000001  # M136: -2274 -> -3728.400000000001
$ws.Cells.Item(136, 14).Value = -8025  # N136: -8059.5 -> -8025

# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(14, 8).Value = 8977.4  # H14: 22349.691 -> 8977.4
$ws.Cells.Item(14, 9).Value = 8977.4  # I14: 22349.691 -> 8977.4
$ws.Cells.Item(14, 11).Value = 26932.2  # K14: 67049.073 -> 26932.2
$ws.Cells.Item(14, 13).Value = -26759.2  # M14: -66876.073 -> -26759.2
$ws.Cells.Item(23, 8).Value = 344.66666  # H23: 531.63635 -> 344.66666
$ws.Cells.Item(23, 9).Value = 489.5  # I23: 799 -> 489.5
$ws.Cells.Item(23, 10).Value = 315.7  # J23: 504.9 -> 315.7
$ws.Cells.Item(23, 11).Value = 1468.5  # K23: 2397 -> 1468.5
$ws.Cells.Item(23, 12).Value = 947.0999999999999  # L23: 1514.7 -> 947.0999999999999
$ws.Cells.Item(23, 13).Value = -1233.5  # M23: -2162 -> -1233.5
$ws.Cells.Item(23, 14).Value = -1417.1  # N23: -1984.7 -> -1417.1
$ws.Cells.Item(39, 8).Value = 3839.2  # H39: 4999.3335 -> 3839.2
$ws.Cells.Item(39, 9).Value = 199  # I39: 0 -> 199
$ws.Cells.Item(39, 10).Value = 4749.25  # J39: 4999.3335 -> 4749.25
$ws.Cells.Item(39, 11).Value = 597  # K39: 0 -> 597
$ws.Cells.Item(39, 12).Value = 14247.75  # L39: 14998.0005 -> 14247.75
$ws.Cells.Item(39, 13).Value = -303  # M39: None -> -303
$ws.Cells.Item(39, 14).Value = -14835.75  # N39: -15586.0005 -> -14835.75
$ws.Cells.Item(61, 8).Value = 302.85715  # H61: 300 -> 302.85715
$ws.Cells.Item(61, 10).Value = 280  # J61: 260 -> 280
$ws.Cells.Item(61, 12).Value = 840  # L61: 780 -> 840
$ws.Cells.Item(61, 14).Value = -1270  # N61: -1210 -> -1270
$ws.Cells.Item(80, 8).Value = 2797  # H80: 2699.2 -> 2797
$ws.Cells.Item(80, 9).Value = 2501  # I80: 2500.6667 -> 2501
$ws.Cells.Item(80, 10).Value = 2994.3333  # J80: 2997 -> 2994.3333
$ws.Cells.Item(80, 11).Value = 7503  # K80: 7502.000100000001 -> 7503
$ws.Cells.Item(80, 12).Value = 8982.999899999999  # L80: 8991 -> 8982.999899999999
$ws.Cells.Item(80, 13).Value = -6567  # M80: -6566.000100000001 -> -6567
$ws.Cells.Item(80, 14).Value = -10854.9999  # N80: -10863 -> -10854.9999
$ws.Cells.Item(83, 8).Value = 2797  # H83: 2699.2 -> 2797
$ws.Cells.Item(83, 9).Value = 2501  # I83: 2500.6667 -> 2501
$ws.Cells.Item(83, 10).Value = 2994.3333  # J83: 2997 -> 2994.3333
$ws.Cells.Item(83, 11).Value = 22509  # K83: 22506.0003 -> 22509
$ws.Cells.Item(83, 12).Value = 26948.9997  # L83: 26973 -> 26948.9997
$ws.Cells.Item(83, 13).Value = -17829  # M83: -17826.0003 -> -17829
$ws.Cells.Item(83, 14).Value = -36308.9997  # N83: -36333 -> -36308.9997
$ws.Cells.Item(132, 8).Value = 1185.75  # H132: 1216 -> 1185.75
$ws.Cells.Item(132, 10).Value = 1097.5  # J132: 1100 -> 1097.5
$ws.Cells.Item(132, 12).Value = 9877.5  # L132: 9900 -> 9877.5
$ws.Cells.Item(132, 14).Value = -14937.5  # N132: -14960 -> -14937.5

# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(57, 8).Value = 0  # H57: 25000 -> 0
$ws.Cells.Item(57, 10).Value = 0  # J57: 25000 -> 0
$ws.Cells.Item(57, 12).Value = 0  # L57: 25000 -> 0
$ws.Cells.Item(57, 14).ClearContents()  # N57: -26640 -> (removed)
$ws.Cells.Item(70, 8).Value = 5399.857  # H70: 6099.8335 -> 5399.857
$ws.Cells.Item(70, 9).Value = 4659.8  # I70: 5524.75 -> 4659.8
$ws.Cells.Item(70, 11).Value = 4659.8  # K70: 5524.75 -> 4659.8
$ws.Cells.Item(70, 13).Value = -4389.8  # M70: -5254.75 -> -4389.8
$ws.Cells.Item(73, 8).Value = 5399.857  # H73: 6099.8335 -> 5399.857
$ws.Cells.Item(73, 9).Value = 4659.8  # I73: 5524.75 -> 4659.8
$ws.Cells.Item(73, 11).Value = 4659.8  # K73: 5524.75 -> 4659.8
$ws.Cells.Item(73, 13).Value = -3723.8  # M73: -4588.75 -> -3723.8
$ws.Cells.Item(80, 8).Value = 9192.111000000001  # H80: 9999.666999999999 -> 9192.111000000001
$ws.Cells.Item(80, 9).Value = 8549.799999999999  # I80: 9999 -> 8549.799999999999
$ws.Cells.Item(80, 10).Value = 9995  # J80: 10000 -> 9995
$ws.Cells.Item(80, 11).Value = 8549.799999999999  # K80: 9999 -> 8549.799999999999
$ws.Cells.Item(80, 12).Value = 9995  # L80: 10000 -> 9995
$ws.Cells.Item(80, 13).Value = -7551.799999999999  # M80: -9001 -> -7551.799999999999
$ws.Cells.Item(80, 14).Value = -11991  # N80: -11996 -> -11991
$ws.Cells.Item(83, 8).Value = 9192.111000000001  # H83: 9999.666999999999 -> 9192.111000000001
$ws.Cells.Item(83, 9).Value = 8549.799999999999  # I83: 9999 -> 8549.799999999999
$ws.Cells.Item(83, 10).Value = 9995  # J83: 10000 -> 9995
$ws.Cells.Item(83, 11).Value = 42749  # K83: 49995 -> 42749
$ws.Cells.Item(83, 12).Value = 49975  # L83: 50000 -> 49975
$ws.Cells.Item(83, 13).Value = -37757  # M83: -45003 -> -37757
$ws.Cells.Item(83, 14).Value = -59959  # N83: -59984 -> -59959
$ws.Cells.Item(102, 8).Value = 2769.6667  # H102: 3063.6 -> 2769.6667
$ws.Cells.Item(102, 10).Value = 1300  # J102: 0 -> 1300
$ws.Cells.Item(102, 12).Value = 1300  # L102: 0 -> 1300
$ws.Cells.Item(102, 14).Value = -4544  # N102: None -> -4544
$ws.Cells.Item(132, 8).Value = 2328.6667  # H132: 2000 -> 2328.6667
$ws.Cells.Item(132, 9).Value = 2328.6667  # I132: 2000 -> 2328.6667
$ws.Cells.Item(132, 11).Value = 6986.000100000001  # K132: 6000 -> 6986.000100000001
$ws.Cells.Item(132, 13).Value = -4456.000100000001  # M132: -3470 -> -4456.000100000001

# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(16, 8).Value = 3899.75  # H16: 3449.6667 -> 3899.75
$ws.Cells.Item(16, 9).Value = 3899.75  # I16: 3449.6667 -> 3899.75
$ws.Cells.Item(16, 11).Value = 3899.75  # K16: 3449.6667 -> 3899.75
$ws.Cells.Item(16, 13).Value = -3729.75  # M16: -3279.6667 -> -3729.75
$ws.Cells.Item(40, 8).Value = 5656  # H40: 5921.846 -> 5656
$ws.Cells.Item(40, 9).Value = 5782.1665  # I40: 6107.8184 -> 5782.1665
$ws.Cells.Item(40, 11).Value = 5782.1665  # K40: 6107.8184 -> 5782.1665
$ws.Cells.Item(40, 13).Value = -5646.1665  # M40: -5971.8184 -> -5646.1665
$ws.Cells.Item(46, 8).Value = 2804.3635  # H46: 2604.3635 -> 2804.3635
$ws.Cells.Item(46, 9).Value = 0  # I46: 300 -> 0
$ws.Cells.Item(46, 10).Value = 2804.3635  # J46: 2834.8 -> 2804.3635
$ws.Cells.Item(46, 11).Value = 0  # K46: 300 -> 0
$ws.Cells.Item(46, 12).Value = 2804.3635  # L46: 2834.8 -> 2804.3635
$ws.Cells.Item(46, 13).ClearContents()  # M46: -112 -> (removed)
$ws.Cells.Item(46, 14).Value = -3180.3635  # N46: -3210.8 -> -3180.3635
$ws.Cells.Item(55, 8).Value = 948.75  # H55: 1560 -> 948.75
$ws.Cells.Item(55, 9).Value = 1270  # I55: 1916.6666 -> 1270
$ws.Cells.Item(55, 10).Value = 413.33334  # J55: 490 -> 413.33334
$ws.Cells.Item(55, 11).Value = 1270  # K55: 1916.6666 -> 1270
$ws.Cells.Item(55, 12).Value = 413.33334  # L55: 490 -> 413.33334
$ws.Cells.Item(55, 13).Value = -1097  # M55: -1743.6666 -> -1097
$ws.Cells.Item(55, 14).Value = -759.33334  # N55: -836 -> -759.33334
$ws.Cells.Item(68, 8).Value = 750  # H68: 1000 -> 750
$ws.Cells.Item(68, 9).Value = 600  # I68: 0 -> 600
$ws.Cells.Item(68, 10).Value = 825  # J68: 1000 -> 825
$ws.Cells.Item(68, 11).Value = 600  # K68: 0 -> 600
$ws.Cells.Item(68, 12).Value = 825  # L68: 1000 -> 825
$ws.Cells.Item(68, 13).Value = 149  # M68: None -> 149
$ws.Cells.Item(68, 14).Value = -2323  # N68: -2498 -> -2323
$ws.Cells.Item(71, 8).Value = 750  # H71: 1000 -> 750
$ws.Cells.Item(71, 9).Value = 600  # I71: 0 -> 600
$ws.Cells.Item(71, 10).Value = 825  # J71: 1000 -> 825
$ws.Cells.Item(71, 11).Value = 3000  # K71: 0 -> 3000
$ws.Cells.Item(71, 12).Value = 4125  # L71: 5000 -> 4125
$ws.Cells.Item(71, 13).Value = 744  # M71: None -> 744
$ws.Cells.Item(71, 14).Value = -11613  # N71: -12488 -> -11613
$ws.Cells.Item(122, 8).Value = 3504  # H122: 0 -> 3504
$ws.Cells.Item(122, 9).Value = 3504  # I122: 0 -> 3504
$ws.Cells.Item(122, 11).Value = 10512  # K122: 0 -> 10512
$ws.Cells.Item(122, 13).Value = -8062  # M122: None -> -8062
$ws.Cells.Item(123, 8).Value = 80000  # H123: 0 -> 80000
$ws.Cells.Item(123, 10).Value = 80000  # J123: 0 -> 80000
$ws.Cells.Item(123, 12).Value = 80000  # L123: 0 -> 80000
$ws.Cells.Item(123, 14).Value = -89800  # N123: None -> -89800
$ws.Cells.Item(132, 8).Value = 16326  # H132: 18081.75 -> 16326
$ws.Cells.Item(132, 9).Value = 19498.375  # I132: 20125.572 -> 19498.375
$ws.Cells.Item(132, 10).Value = 3636.5  # J132: 3775 -> 3636.5
$ws.Cells.Item(132, 11).Value = 58495.125  # K132: 60376.716 -> 58495.125
$ws.Cells.Item(132, 12).Value = 10909.5  # L132: 11325 -> 10909.5
$ws.Cells.Item(132, 13).Value = -55965.125  # M132: -57846.716 -> -55965.125
$ws.Cells.Item(132, 14).Value = -15969.5  # N132: -16385 -> -15969.5
$ws.Cells.Item(136, 8).Value = 3408.875  # H136: 3630.8572 -> 3408.875
$ws.Cells.Item(136, 9).Value = 3227.6667  # I136: 3502.2 -> 3227.6667
$ws.Cells.Item(136, 11).Value = 9683.000100000001  # K136: 10506.6 -> 9683.000100000001
$ws.Cells.Item(136, 13).Value = -7133.000100000001  # M136: -7956.599999999999 -> -7133.000100000001

# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(80, 8).Value = 33650.5  # H80: 0 -> 33650.5
$ws.Cells.Item(80, 10).Value = 33650.5  # J80: 0 -> 33650.5
$ws.Cells.Item(80, 12).Value = 33650.5  # L80: 0 -> 33650.5
$ws.Cells.Item(80, 14).Value = -35646.5  # N80: None -> -35646.5
$ws.Cells.Item(83, 8).Value = 33650.5  # H83: 0 -> 33650.5
$ws.Cells.Item(83, 10).Value = 33650.5  # J83: 0 -> 33650.5
$ws.Cells.Item(83, 12).Value = 100951.5  # L83: 0 -> 100951.5
$ws.Cells.Item(83, 14).Value = -110935.5  # N83: None -> -110935.5
$ws.Cells.Item(132, 8).Value = 1579.6  # H132: 1699.3334 -> 1579.6
$ws.Cells.Item(132, 9).Value = 1474.75  # I132: 1549.5 -> 1474.75
$ws.Cells.Item(132, 11).Value = 4424.25  # K132: 4648.5 -> 4424.25
$ws.Cells.Item(132, 13).Value = -1894.25  # M132: -2118.5 -> -1894.25
